# Updates the cryptos list values to match the latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.246.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +9.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.259.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.55%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "394.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("E6").Value = "  +5.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.254.57"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.44%  "

$ws.Range("E8").Value = "  +5.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.623"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "39.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0968"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.56%  "

$ws.Range("E13").Value = "  +2.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.768.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.266.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.56%  "

$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "57.027.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.39%  "

$ws.Range("E22").Value = "  +8.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "301.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.06"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.01%  "

$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.99%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").Value = "  +1.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "37.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0484"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +15.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "134.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.03%  "

$ws.Range("E43").Value = "  +2.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.03%  "

$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.120"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.285"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.150.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +37.01%  "
